# "statusbericht stunden eingetragen und to do notizen von dalessandro"
#
# Enters the latest weekly hours (week "W19.3" = column R, corrected; and
# week "W19.5" = column T, newly reported) for the team members' booked
# hours on the "Kosten" sheet, for both the "4. Realisierung" (rows 67-72)
# and "5. Validierung" (rows 85-90) work packages. All dependent SUM /
# running-total formulas (AA, row 73/74/75/76, row 91/92/93/94) recompute
# automatically on recalculation.
#
# It also moves the active selection/tab the way the author left the file:
# selecting T73 on "Kosten" and ending with "Status-4" as the active sheet.

$wb = $excel.ActiveWorkbook
$kosten = $wb.Worksheets.Item("Kosten")

# --- Realisierung (rows 67-72): corrections to W19.3 (R) and new entries
#     for W19.5 (T) ---------------------------------------------------------
$kosten.Range("R67").Value = 10
$kosten.Range("T67").Value = 9

$kosten.Range("T68").Value = 10

$kosten.Range("R69").Value = 9
$kosten.Range("T69").Value = 1

$kosten.Range("R70").Value = 9
$kosten.Range("T70").Value = 9

$kosten.Range("R71").Value = 9
$kosten.Range("T71").Value = 1

$kosten.Range("R72").Value = 9
$kosten.Range("T72").Value = 9

# --- Validierung (rows 85-90): new entries for W19.5 (T) -------------------
$kosten.Range("T86").Value = 9
$kosten.Range("T87").Value = 9
$kosten.Range("T89").Value = 9

# Recalculate so every dependent SUM / running-total formula picks up the
# new figures.
$excel.Calculate()

# --- Leave the view the way it ended up after the edits --------------------
$kosten.Activate()
$kosten.Range("T73").Select()

$status4 = $wb.Worksheets.Item("Status-4")
$status4.Activate()
